$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row 1 cells: "_old" suffix -> "_FV2304", "_new" suffix -> "_FV2310"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2304"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2310"
        }
    }
}

# Turn the used range into an Excel table ("Table1")
$range = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Application.ActiveWindow.FreezePanes = $true
